# Replace the text content of a run inside a paragraph while preserving the
# paragraph's exact run/paragraph-property structure (including any "empty"
# <w:r/> runs that a plain Find.Execute text replace would otherwise merge
# away). Every paragraph whose visible text contains $oldText has that exact
# run-text swapped for $newText.
function Replace-ParagraphText {
    param([string]$oldText, [string]$newText)

    $d = $word.ActiveDocument
    $oldRunText = "<w:t>" + $oldText + "</w:t>"
    $oldRunTextPreserve = '<w:t xml:space="preserve">' + $oldText + "</w:t>"
    $replacedAny = $false

    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $full = $p.Range

        if ($full.Text.Contains($oldText)) {
            $x = $full.WordOpenXML

            $idx1 = $x.IndexOf("<w:p ")
            if ($idx1 -lt 0) { $idx1 = $x.IndexOf("<w:p>") }
            $idx2 = $x.IndexOf("</w:p>") + 6
            $pXml = $x.Substring($idx1, $idx2 - $idx1)

            # Normalize the opening <w:p ...> tag back down to a bare <w:p>
            # (the WordOpenXML getter stamps synthetic paraId/rsid
            # attributes on it that the source document never had).
            $closeIdx = $pXml.IndexOf(">")
            $pXml = "<w:p>" + $pXml.Substring($closeIdx + 1)

            if ($pXml.Contains($oldRunText)) {
                $newPXml = $pXml.Replace($oldRunText, "<w:t>" + $newText + "</w:t>")
            } elseif ($pXml.Contains($oldRunTextPreserve)) {
                $newPXml = $pXml.Replace($oldRunTextPreserve, '<w:t xml:space="preserve">' + $newText + "</w:t>")
            } else {
                Write-Output ("ERROR: run text not found verbatim in paragraph " + $i)
                continue
            }

            $wrapped = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + "<w:body>" + $newPXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

            $countBefore = $d.Paragraphs.Count

            $collapsed = $full.Duplicate
            $collapsed.Collapse(0)
            $collapsed.InsertXML($wrapped)

            # InsertXML-ing a whole <w:p> at the very end of the document
            # body (there is no following paragraph to reuse as the
            # insertion boundary) spawns one extra trailing empty
            # paragraph instead of replacing in place. Collapse it back
            # out by deleting the paragraph-mark boundary that separates
            # our freshly written paragraph from the stray empty one(s).
            $countAfter = $d.Paragraphs.Count
            while ($countAfter -gt $countBefore) {
                $fixedP = $d.Paragraphs.Item($countAfter - 1)
                $strayP = $d.Paragraphs.Item($countAfter)
                $boundary = $d.Range($fixedP.Range.End - 1, $strayP.Range.End)
                $boundary.Delete()
                $countAfter = $d.Paragraphs.Count
            }

            $replacedAny = $true
        }
    }

    if (-not $replacedAny) {
        Write-Output ("ERROR: old text not found anywhere: " + $oldText)
    }
}

# 1. Title heading + bold title near the end (same old->new text; replaces both occurrences)
Replace-ParagraphText "Play Lightning Horseman Slot for Free - Review" "Play Lightning Horseman Slot Game Free"

# "What we like" bullet list
Replace-ParagraphText "Exciting Lightning Respin feature" "Thrilling gameplay based on the classic tale of Sleepy Hollow"
Replace-ParagraphText "High-paying headless rider wild symbol" "High-paying wild symbol with potential for big jackpots"
Replace-ParagraphText "Big jackpots during bonus spins" "Exciting Lightning Respin feature with locked symbols and potential prizes"
Replace-ParagraphText "Produced by respected developer Lightning Box Games" "Eerie graphics and atmospheric soundtrack"

# "What we don't like" bullet list
Replace-ParagraphText "Graphics may feel dated to some players" "Graphics may seem dated to some players"
Replace-ParagraphText "Bonus features may be hard to trigger" "Limited number of paylines"

# Meta description paragraph
Replace-ParagraphText "Learn all about Lightning Horseman, the gothic-themed slot game by Lightning Box Games, and play it for free. Discover the game's interesting features and big jackpots." "Read our review of Lightning Horseman slot game and play for free. Experience thrilling gameplay and big jackpots."
